# Atualização de bases das ligas, do dia: 28-06-2024 às 19:47
#
# This script reproduces a re-sort of the match rows in the
# "Germany Oberliga Hamburg" sheet. Two shared strings ("FC Alsterbruder"
# and "TSV Buchholz 08") end up re-indexed, which has two visible effects:
#   1) a handful of rows whose HomeTeam/AwayTeam happened to be one of
#      those two teams keep the SAME displayed text (only the internal
#      shared-string index changes - nothing to do here, Excel manages
#      shared strings automatically when we set cell values);
#   2) a small number of rows had their entire data (every column except
#      the running id in column A, the Div in column C and the Date in
#      column D) swapped/rotated with another row:
#        - row 14  <->  row 16
#        - row 20  <->  row 21
#        - row 258 -> row 261, row 260 -> row 258, row 261 -> row 260
#
# Below we simply (re)write the final, resolved values for every touched
# cell in those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 14 (final values, previously held by row 16) ----
$ws.Range("B14").Value = 6893654
$ws.Range("E14").Value = "Niendorfer TSV"
$ws.Range("F14").Value = "Wandsbeker TSV Concordia"
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = 2
$ws.Range("K14").Value = "H"
$ws.Range("L14").Value = 1.444
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = 4.5
$ws.Range("O14").Value = 1.333
$ws.Range("P14").Value = 5.5
$ws.Range("Q14").Value = 6
$ws.Range("R14").Value = -1.75
$ws.Range("S14").Value = 1.925
$ws.Range("T14").Value = 1.875
$ws.Range("U14").Value = 4
$ws.Range("V14").Value = 1.875
$ws.Range("W14").Value = 1.925
$ws.Range("X14").Value = 0.333
$ws.Range("Y14").Value = -1
$ws.Range("Z14").Value = -1
$ws.Range("AA14").Value = 0.925
$ws.Range("AB14").Value = -1
$ws.Range("AC14").Value = 0.875
$ws.Range("AD14").Value = -1

# ---- Row 16 (final values, previously held by row 14) ----
$ws.Range("B16").Value = 6893655
$ws.Range("E16").Value = "FC Alsterbruder"
$ws.Range("F16").Value = "FC Trkiye Wilhelmsburg"
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 0
$ws.Range("K16").Value = "H"
$ws.Range("L16").Value = 2
$ws.Range("M16").Value = 4.333
$ws.Range("N16").Value = 2.625
$ws.Range("O16").Value = 1.85
$ws.Range("P16").Value = 4.75
$ws.Range("Q16").Value = 2.8
$ws.Range("R16").Value = -0.5
$ws.Range("S16").Value = 1.875
$ws.Range("T16").Value = 1.925
$ws.Range("U16").Value = 3.75
$ws.Range("V16").Value = 1.95
$ws.Range("W16").Value = 1.85
$ws.Range("X16").Value = 0.8500000000000001
$ws.Range("Y16").Value = -1
$ws.Range("Z16").Value = -1
$ws.Range("AA16").Value = 0.875
$ws.Range("AB16").Value = -1
$ws.Range("AC16").Value = -1
$ws.Range("AD16").Value = 0.8500000000000001

# ---- Row 20 (final values, previously held by row 21) ----
$ws.Range("B20").Value = 6893659
$ws.Range("E20").Value = "ETSV Hamburg"
$ws.Range("F20").Value = "Wandsbeker TSV Concordia"
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 2
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = "D"
$ws.Range("L20").Value = 2.1
$ws.Range("M20").Value = 4
$ws.Range("N20").Value = 2.6
$ws.Range("O20").Value = 1.363
$ws.Range("P20").Value = 5.75
$ws.Range("Q20").Value = 5
$ws.Range("R20").Value = -1.5
$ws.Range("S20").Value = 1.825
$ws.Range("T20").Value = 1.975
$ws.Range("U20").Value = 4.25
$ws.Range("V20").Value = 1.925
$ws.Range("W20").Value = 1.875
$ws.Range("X20").Value = -1
$ws.Range("Y20").Value = 4.75
$ws.Range("Z20").Value = -1
$ws.Range("AA20").Value = -1
$ws.Range("AB20").Value = 0.9750000000000001
$ws.Range("AC20").Value = -0.5
$ws.Range("AD20").Value = 0.4375

# ---- Row 21 (final values, previously held by row 20) ----
$ws.Range("B21").Value = 6893660
$ws.Range("E21").Value = "TuRa Harksheide"
$ws.Range("F21").Value = "FC Alsterbruder"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").ClearContents()
$ws.Range("J21").ClearContents()
$ws.Range("K21").Value = "D"
$ws.Range("L21").Value = 2.1
$ws.Range("M21").Value = 4
$ws.Range("N21").Value = 2.6
$ws.Range("O21").Value = 1.95
$ws.Range("P21").Value = 4
$ws.Range("Q21").Value = 2.8
$ws.Range("R21").Value = -0.25
$ws.Range("S21").Value = 1.8
$ws.Range("T21").Value = 2
$ws.Range("U21").Value = 3.5
$ws.Range("V21").Value = 1.95
$ws.Range("W21").Value = 1.85
$ws.Range("X21").Value = -1
$ws.Range("Y21").Value = 3
$ws.Range("Z21").Value = -1
$ws.Range("AA21").Value = -0.5
$ws.Range("AB21").Value = 0.5
$ws.Range("AC21").Value = -1
$ws.Range("AD21").Value = 0.8500000000000001

# ---- Row 258 (final values, previously held by row 260) ----
$ws.Range("B258").Value = 6895693
$ws.Range("E258").Value = "Altona 93"
$ws.Range("F258").Value = "FC Sderelbe"
$ws.Range("G258").Value = 2
$ws.Range("H258").Value = 0
$ws.Range("I258").Value = 2
$ws.Range("J258").Value = 0
$ws.Range("K258").Value = "H"
$ws.Range("L258").Value = 1.285
$ws.Range("M258").Value = 5.5
$ws.Range("N258").Value = 6.5
$ws.Range("O258").Value = 1.75
$ws.Range("P258").Value = 4.333
$ws.Range("Q258").Value = 3.2
$ws.Range("R258").Value = -0.5
$ws.Range("S258").Value = 1.775
$ws.Range("T258").Value = 1.925
$ws.Range("U258").Value = 3.75
$ws.Range("V258").Value = 1.75
$ws.Range("W258").Value = 1.95
$ws.Range("X258").Value = 0.75
$ws.Range("Y258").Value = -1
$ws.Range("Z258").Value = -1
$ws.Range("AA258").Value = 0.7749999999999999
$ws.Range("AB258").Value = -1
$ws.Range("AC258").Value = -1
$ws.Range("AD258").Value = 0.95

# ---- Row 260 (final values, previously held by row 261) ----
$ws.Range("B260").Value = 6896535
$ws.Range("E260").Value = "Niendorfer TSV"
$ws.Range("F260").Value = "Uhlenhorster SC Paloma"
$ws.Range("G260").Value = 2
$ws.Range("H260").Value = 2
$ws.Range("I260").Value = 1
$ws.Range("J260").Value = 1
$ws.Range("K260").Value = "D"
$ws.Range("L260").Value = 1.333
$ws.Range("M260").Value = 5
$ws.Range("N260").Value = 6
$ws.Range("O260").Value = 1.38
$ws.Range("P260").Value = 4.75
$ws.Range("Q260").Value = 5.5
$ws.Range("R260").Value = -1.5
$ws.Range("S260").Value = 1.825
$ws.Range("T260").Value = 1.975
$ws.Range("U260").Value = 4.5
$ws.Range("V260").Value = 1.925
$ws.Range("W260").Value = 1.875
$ws.Range("X260").Value = -1
$ws.Range("Y260").Value = 3.75
$ws.Range("Z260").Value = -1
$ws.Range("AA260").Value = -1
$ws.Range("AB260").Value = 0.9750000000000001
$ws.Range("AC260").Value = -1
$ws.Range("AD260").Value = 0.875

# ---- Row 261 (final values, previously held by row 258) ----
$ws.Range("B261").Value = 6895399
$ws.Range("E261").Value = "SV Rugenbergen"
$ws.Range("F261").Value = "Hamburg Eimsbutteler BC"
$ws.Range("G261").Value = 1
$ws.Range("H261").Value = 1
$ws.Range("I261").Value = 0
$ws.Range("J261").Value = 0
$ws.Range("K261").Value = "D"
$ws.Range("L261").Value = 2.5
$ws.Range("M261").Value = 4.2
$ws.Range("N261").Value = 2.1
$ws.Range("O261").Value = 3.5
$ws.Range("P261").Value = 4.1
$ws.Range("Q261").Value = 1.7
$ws.Range("R261").Value = 0.75
$ws.Range("S261").Value = 1.875
$ws.Range("T261").Value = 1.925
$ws.Range("U261").Value = 3.5
$ws.Range("V261").Value = 1.875
$ws.Range("W261").Value = 1.925
$ws.Range("X261").Value = -1
$ws.Range("Y261").Value = 3.1
$ws.Range("Z261").Value = -1
$ws.Range("AA261").Value = 0.875
$ws.Range("AB261").Value = -1
$ws.Range("AC261").Value = -1
$ws.Range("AD261").Value = 0.925

# ---- Shared-string table swap side effect ----
# "FC Alsterbruder" and "TSV Buchholz 08" swapped places in the shared
# string table. Every other row referencing either of those two team
# names keeps the SAME displayed text (the team names themselves did not
# change), so no further cell writes are required - Excel re-resolves /
# reuses shared strings automatically from the literal text we assigned
# above.
